# "Fruta / hortaliza, semanal"
#
# A new weekly price record for "Macroferia Regional de Talca - Haba" needs
# to be inserted as the data row right after the header/first record, i.e.
# at worksheet row 34 (the table's data rows start at row 2). Inserting a
# row there pushes the former rows 34..103 down to 35..104, growing the
# used range from A1:R103 to A1:R104, and the brand-new row 34 is then
# populated with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (old rows 34-103) down to rows 35-104.
$ws.Rows(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 44868
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 100112026
$ws.Range("G34").Value = "Haba"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 7000
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Región del Maule"
$ws.Range("P34").Value = 280
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
